# Weekly price update: insert a new data row for "Terminal Hortofrutícola
# Agro Chillán - Piña" ahead of the existing rows (new week, 2021-11-09),
# shifting the previously-existing rows 168-172 down to 169-173.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 168..172 down to 169..173, leaving a blank row 168 to fill in.
$ws.Rows("168:168").Insert()

$ws.Cells.Item(168, 1).Value = 7
$ws.Cells.Item(168, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(168, 3).Value = "Ñuble"
$ws.Cells.Item(168, 4).Value = 44509
$ws.Cells.Item(168, 5).Value = 16
$ws.Cells.Item(168, 6).Value = "Fruta"
$ws.Cells.Item(168, 7).Value = 100108
$ws.Cells.Item(168, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(168, 9).Value = 100108005
$ws.Cells.Item(168, 10).Value = "Piña"
$ws.Cells.Item(168, 11).Value = "Caramelo"
$ws.Cells.Item(168, 12).Value = "Segunda"
$ws.Cells.Item(168, 13).Value = 120
$ws.Cells.Item(168, 14).Value = 18000
$ws.Cells.Item(168, 15).Value = 19000
$ws.Cells.Item(168, 16).Value = 18500
$ws.Cells.Item(168, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(168, 18).Value = "Ecuador"
$ws.Cells.Item(168, 19).Value = 1321
$ws.Cells.Item(168, 20).Value = 14
